# Update latest optimisation output (run 152)
$wb = $excel.ActiveWorkbook

# --- Schedule sheet ---
$ws1 = $wb.Worksheets.Item("Schedule")
$ws1.Range("B2").Value = 46043.20833333334
$ws1.Range("C2").Value = 5
$ws1.Range("D2").Value = 18.9
$ws1.Range("E2").Value = 645.2629852499999
$ws1.Range("F2").Value = 34.14089869047619
$ws1.Range("A4").Value = 46043.91666666666
$ws1.Range("C4").Value = 4.5
$ws1.Range("D4").Value = 17.01
$ws1.Range("E4").Value = 566.1105060000001
$ws1.Range("F4").Value = 33.2810409171076
$ws1.Range("E5").Value = 20.79376649999999
$ws1.Range("F5").Value = 0.5790522556390975

# --- Detailed sheet ---
$ws2 = $wb.Worksheets.Item("Detailed")
$ws2.Range("E11").Value = "ON"
$ws2.Range("B39").Value = 61.05003
$ws2.Range("B40").Value = 79.95
$ws2.Range("B41").Value = 76.60739
$ws2.Range("C41").Value = "historical"
$ws2.Range("B42").Value = 79.95003
$ws2.Range("C42").Value = "historical"
$ws2.Range("B43").Value = 77.94
$ws2.Range("B44").Value = 79.95
$ws2.Range("B45").Value = 77.94
$ws2.Range("E45").Value = "OFF"
$ws2.Range("B46").Value = 65.51413
$ws2.Range("B47").Value = 64.8901
$ws2.Range("B48").Value = 64.8901
$ws2.Range("B49").Value = 64.8901
$ws2.Range("B50").Value = 64.8901
$ws2.Range("B51").Value = 65.64485
$ws2.Range("B52").Value = 64.50382
$ws2.Range("B53").Value = 62.6453
$ws2.Range("B54").Value = 62.75766
$ws2.Range("B55").Value = 63.6179
$ws2.Range("B56").Value = 66.58355
$ws2.Range("B57").Value = 66.18547
$ws2.Range("B58").Value = 66.5827
$ws2.Range("B59").Value = 67.75791
$ws2.Range("B60").Value = 73.2
$ws2.Range("B61").Value = 76.74068
$ws2.Range("B64").Value = 35.88
$ws2.Range("B66").Value = -5.74313
$ws2.Range("B67").Value = -7.69633
$ws2.Range("B68").Value = -10.11673
$ws2.Range("B69").Value = -13.49924
$ws2.Range("B70").Value = -7.94851
$ws2.Range("B71").Value = -9.86476
$ws2.Range("B72").Value = -8.91861
$ws2.Range("B73").Value = -6.48882
$ws2.Range("B74").Value = -9.65649
$ws2.Range("B75").Value = -7.9049
$ws2.Range("B76").Value = -9.40225
$ws2.Range("B77").Value = -5.95261
$ws2.Range("B78").Value = -5.50985
$ws2.Range("B79").Value = 0.00917
$ws2.Range("B83").Value = -10
$ws2.Range("B84").Value = -11.84289
$ws2.Range("B85").Value = -8.08329
$ws2.Range("B86").Value = -0.4592
$ws2.Range("B87").Value = 0.00037
$ws2.Range("B88").Value = 10.40979
$ws2.Range("B90").Value = 53.90468
$ws2.Range("B91").Value = 54.5312
$ws2.Range("B93").Value = 57.01318
$ws2.Range("B94").Value = 42.98129
$ws2.Range("B95").Value = 56.39409
$ws2.Range("B96").Value = 56.20379
$ws2.Range("B97").Value = 48.29198

Write-Host "Applied run 152 updates"
